# [NEW FEATURE] support to set the specified color of word
#
# Adds a third "color" column to the wordcloud data sheet so each word can
# carry an explicit display color alongside its weight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing headers and add the new "color" header.
$ws.Range("A1").Value = "word"
$ws.Range("B1").Value = "weight"
$ws.Range("C1").Value = "color"

# Populate the new color column for the top keywords.
$ws.Range("C2").Value = "#CCFFE5"
$ws.Range("C3").Value = "#9999FF"
$ws.Range("C4").Value = "#000000"

# Leave the selection on the last-edited cell, matching the authored change.
$ws.Range("C4").Select()
